$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column CD holds the "2-nov" data, mirroring the existing "1-nov" column (CC)
$ws.Range("CD1").Value = "2-nov"

$ws.Range("CD2").Value = 10
$ws.Range("CD3").Value = 10
$ws.Range("CD4").Value = 9
$ws.Range("CD5").Value = 5
$ws.Range("CD6").Value = 10
$ws.Range("CD7").Value = 6
$ws.Range("CD8").Value = 12
$ws.Range("CD9").Value = 14
$ws.Range("CD10").Value = 5
$ws.Range("CD11").Value = 0

# Match formatting of the adjacent "1-nov" column (CC) - set alignment
# before number format so the style dedupes onto the existing cellXf
# instead of leaving an orphaned intermediate style behind.
$ws.Range("CD1").NumberFormat = $ws.Range("CC1").NumberFormat

$ws.Range("CD2:CD11").HorizontalAlignment = -4108
$ws.Range("CD2:CD11").NumberFormat = $ws.Range("CC2:CC11").NumberFormat

$ws.Range("CD12").Select()
